$d = $word.ActiveDocument

# Mapping of old text -> new text, applied via Find/Replace on the whole document content.
# Each old text value is unique within the document, so each replacement is unambiguous.
$replacements = @(
    @("2024-06-08 Saturday", "2024-06-09 Sunday"),
    @("914×9=8226", "898×8=7184"),
    @("509×7=3563", "611×3=1833"),
    @("281×4=1124", "338×8=2704"),
    @("822×2=1644", "704×3=2112"),
    @("338×7=2366", "823×5=4115"),
    @("112×4=448", "513×5=2565"),
    @("352×7=2464", "652×7=4564"),
    @("207×5=1035", "196×5=980"),
    @("631×7=4417", "712×8=5696"),
    @("696×2=1392", "995×9=8955"),
    @("840×2=1680", "881×5=4405"),
    @("133×7=931", "912×6=5472"),
    @("641×8=5128", "175×6=1050"),
    @("131×4=524", "514×2=1028"),
    @("385×2=770", "184×2=368"),
    @("990×4=3960", "770×3=2310"),
    @("777×9=6993", "817×7=5719"),
    @("567×4=2268", "349×8=2792"),
    @("755×5=3775", "180×2=360"),
    @("675×9=6075", "245×3=735"),
    @("730×3=2190", "904×4=3616"),
    @("472×2=944", "427×8=3416"),
    @("961×5=4805", "569×9=5121"),
    @("124×7=868", "785×9=7065"),
    @("686×6=4116", "454×3=1362"),
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
